$wb = $excel.ActiveWorkbook

# Sheet names to update: 展览 (Exhibitions) and 全部类型 (All Types)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2302
    $ws.Range("F3").Value = 1752
    $ws.Range("F4").Value = 341
    $ws.Range("F6").Value = 905
    $ws.Range("F8").Value = 5862
}
